$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A1").ClearContents()

$new = $wb.Worksheets.Add($null, $ws1)
$new.Name = "Sheet1"

$src = $ws1.Range("A9:K40")
$dst = $new.Range("A1:K32")
$src.Copy($dst)

$ws1.Range("A9:K40").Select()
